$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-20 18:48:46"
$ws.Range("E3").Value = "2026-02-20 18:48:48"
$ws.Range("O3").Value = "-5.1 °C"
$ws.Range("E4").Value = "2026-02-20 18:48:51"
$ws.Range("J4").Value = "1022.2 hPa"
$ws.Range("E5").Value = "2026-02-20 18:48:54"
$ws.Range("N5").Value = "-6.0 °C 18:27 TU"
$ws.Range("O5").Value = "-4.6 °C"
$ws.Range("E6").Value = "2026-02-20 18:48:56"
$ws.Range("J6").Value = "1022.2 hPa"
$ws.Range("E7").Value = "2026-02-20 18:48:59"
$ws.Range("H7").Value = "45%"
$ws.Range("J7").Value = "1022.1 hPa"
$ws.Range("O7").Value = "13.5 °C"
$ws.Range("E8").Value = "2026-02-20 18:49:01"
$ws.Range("J8").Value = "1022.4 hPa"
$ws.Range("E9").Value = "2026-02-20 18:49:04"
$ws.Range("E10").Value = "2026-02-20 18:49:07"
$ws.Range("E11").Value = "2026-02-20 18:49:09"
$ws.Range("E12").Value = "2026-02-20 18:49:12"
$ws.Range("E13").Value = "2026-02-20 18:49:14"
$ws.Range("J13").Value = "1023.1 hPa"
$ws.Range("E14").Value = "2026-02-20 18:49:17"
$ws.Range("H14").Value = "54%"
$ws.Range("N14").Value = "9.0 °C 18:27 TU"
$ws.Range("E15").Value = "2026-02-20 18:49:19"
$ws.Range("H15").Value = "45%"
$ws.Range("O15").Value = "13.4 °C"
$ws.Range("E16").Value = "2026-02-20 18:49:22"
$ws.Range("M16").Value = "0.6 °C 18:27 TU"
$ws.Range("O16").Value = "-3.5 °C"
$ws.Range("E17").Value = "2026-02-20 18:49:25"
$ws.Range("O17").Value = "2.8 °C"
$ws.Range("E18").Value = "2026-02-20 18:49:28"
$ws.Range("J18").Value = "1022.5 hPa"
$ws.Range("E19").Value = "2026-02-20 18:49:30"
$ws.Range("E20").Value = "2026-02-20 18:49:33"
$ws.Range("O20").Value = "-3.0 °C"
$ws.Range("E21").Value = "2026-02-20 18:49:36"
$ws.Range("J21").Value = "1022.1 hPa"
$ws.Range("O21").Value = "9.6 °C"
$ws.Range("E22").Value = "2026-02-20 18:49:38"
$ws.Range("E23").Value = "2026-02-20 18:49:41"
$ws.Range("H23").Value = "67%"
$ws.Range("M23").Value = "-0.4 °C 18:29 TU"
$ws.Range("O23").Value = "-5.2 °C"
$ws.Range("E24").Value = "2026-02-20 18:49:44"
$ws.Range("J24").Value = "1025.1 hPa"
$ws.Range("E25").Value = "2026-02-20 18:49:47"
$ws.Range("H25").Value = "50%"
$ws.Range("E26").Value = "2026-02-20 18:49:49"
$ws.Range("H26").Value = "33%"
$ws.Range("J26").Value = "1021.3 hPa"
$ws.Range("O26").Value = "5.5 °C"
$ws.Range("E27").Value = "2026-02-20 18:49:52"
$ws.Range("E28").Value = "2026-02-20 18:49:55"
$ws.Range("H28").Value = "63%"
$ws.Range("J28").Value = "1022.5 hPa"
$ws.Range("E29").Value = "2026-02-20 18:49:57"
$ws.Range("O29").Value = "9.5 °C"
$ws.Range("E30").Value = "2026-02-20 18:50:00"
$ws.Range("H30").Value = "57%"
$ws.Range("J30").Value = "1021.9 hPa"
$ws.Range("E31").Value = "2026-02-20 18:50:03"
$ws.Range("J31").Value = "1021.0 hPa"
$ws.Range("E32").Value = "2026-02-20 18:50:06"
$ws.Range("O32").Value = "4.6 °C"
$ws.Range("E33").Value = "2026-02-20 18:50:08"
$ws.Range("H33").Value = "39%"
$ws.Range("J33").Value = "1022.4 hPa"
$ws.Range("E34").Value = "2026-02-20 18:50:11"
$ws.Range("O34").Value = "0.5 °C"
$ws.Range("E35").Value = "2026-02-20 18:50:14"
$ws.Range("J35").Value = "1026.5 hPa"
$ws.Range("E36").Value = "2026-02-20 18:50:16"
$ws.Range("J36").Value = "1022.1 hPa"
$ws.Range("E37").Value = "2026-02-20 18:50:19"
$ws.Range("H37").Value = "64%"
$ws.Range("J37").Value = "1023.9 hPa"
$ws.Range("E38").Value = "2026-02-20 18:50:22"
$ws.Range("E39").Value = "2026-02-20 18:50:24"
$ws.Range("E40").Value = "2026-02-20 18:50:26"
$ws.Range("J40").Value = "1023.0 hPa"
$ws.Range("O40").Value = "10.6 °C"
$ws.Range("E41").Value = "2026-02-20 18:50:29"
$ws.Range("H41").Value = "48%"
$ws.Range("E42").Value = "2026-02-20 18:50:32"
$ws.Range("O42").Value = "10.1 °C"
$ws.Range("E43").Value = "2026-02-20 18:50:34"
$ws.Range("O43").Value = "5.1 °C"
$ws.Range("E44").Value = "2026-02-20 18:50:37"
$ws.Range("H44").Value = "80%"
$ws.Range("M44").Value = "-1.6 °C 18:25 TU"
$ws.Range("O44").Value = "-5.0 °C"
$ws.Range("E45").Value = "2026-02-20 18:50:39"
$ws.Range("J45").Value = "1029.2 hPa"
$ws.Range("E46").Value = "2026-02-20 18:50:42"
